# Updates the 100 arithmetic-equation answer cells in the single table of
# the document, replacing each old "a op b = c" string with the new
# generated equation string (per the commit's regenerated answer set).
$d = $word.ActiveDocument

$d.Content.Find.Execute("88-77=11", $true, $false, $false, $false, $false, $true, 1, $false, "88-70=18", 2) | Out-Null
$d.Content.Find.Execute("40+24=64", $true, $false, $false, $false, $false, $true, 1, $false, "85-80=5", 2) | Out-Null
$d.Content.Find.Execute("46-0=46", $true, $false, $false, $false, $false, $true, 1, $false, "13+85=98", 2) | Out-Null
$d.Content.Find.Execute("9-4=5", $true, $false, $false, $false, $false, $true, 1, $false, "69-38=31", 2) | Out-Null
$d.Content.Find.Execute("64-39=25", $true, $false, $false, $false, $false, $true, 1, $false, "67-9=58", 2) | Out-Null
$d.Content.Find.Execute("90-84=6", $true, $false, $false, $false, $false, $true, 1, $false, "41-21=20", 2) | Out-Null
$d.Content.Find.Execute("13+31=44", $true, $false, $false, $false, $false, $true, 1, $false, "29-29=0", 2) | Out-Null
$d.Content.Find.Execute("37+19=56", $true, $false, $false, $false, $false, $true, 1, $false, "28+48=76", 2) | Out-Null
$d.Content.Find.Execute("70-48=22", $true, $false, $false, $false, $false, $true, 1, $false, "9+44=53", 2) | Out-Null
$d.Content.Find.Execute("9+25=34", $true, $false, $false, $false, $false, $true, 1, $false, "0+48=48", 2) | Out-Null
$d.Content.Find.Execute("95-37=58", $true, $false, $false, $false, $false, $true, 1, $false, "66-64=2", 2) | Out-Null
$d.Content.Find.Execute("11+39=50", $true, $false, $false, $false, $false, $true, 1, $false, "3+34=37", 2) | Out-Null
$d.Content.Find.Execute("95-88=7", $true, $false, $false, $false, $false, $true, 1, $false, "27+26=53", 2) | Out-Null
$d.Content.Find.Execute("10+50=60", $true, $false, $false, $false, $false, $true, 1, $false, "38-25=13", 2) | Out-Null
$d.Content.Find.Execute("21+58=79", $true, $false, $false, $false, $false, $true, 1, $false, "2+31=33", 2) | Out-Null
$d.Content.Find.Execute("29+44=73", $true, $false, $false, $false, $false, $true, 1, $false, "13+60=73", 2) | Out-Null
$d.Content.Find.Execute("78-52=26", $true, $false, $false, $false, $false, $true, 1, $false, "38-31=7", 2) | Out-Null
$d.Content.Find.Execute("39+53=92", $true, $false, $false, $false, $false, $true, 1, $false, "94-85=9", 2) | Out-Null
$d.Content.Find.Execute("77-77=0", $true, $false, $false, $false, $false, $true, 1, $false, "29-11=18", 2) | Out-Null
$d.Content.Find.Execute("34+49=83", $true, $false, $false, $false, $false, $true, 1, $false, "90-30=60", 2) | Out-Null
$d.Content.Find.Execute("45+21=66", $true, $false, $false, $false, $false, $true, 1, $false, "57-9=48", 2) | Out-Null
$d.Content.Find.Execute("22+77=99", $true, $false, $false, $false, $false, $true, 1, $false, "56-50=6", 2) | Out-Null
$d.Content.Find.Execute("33+63=96", $true, $false, $false, $false, $false, $true, 1, $false, "97-33=64", 2) | Out-Null
$d.Content.Find.Execute("93-34=59", $true, $false, $false, $false, $false, $true, 1, $false, "45-9=36", 2) | Out-Null
$d.Content.Find.Execute("32+62=94", $true, $false, $false, $false, $false, $true, 1, $false, "47+48=95", 2) | Out-Null
$d.Content.Find.Execute("70+29=99", $true, $false, $false, $false, $false, $true, 1, $false, "93-12=81", 2) | Out-Null
$d.Content.Find.Execute("19+16=35", $true, $false, $false, $false, $false, $true, 1, $false, "50-23=27", 2) | Out-Null
$d.Content.Find.Execute("3+70=73", $true, $false, $false, $false, $false, $true, 1, $false, "51-8=43", 2) | Out-Null
$d.Content.Find.Execute("66-0=66", $true, $false, $false, $false, $false, $true, 1, $false, "38+1=39", 2) | Out-Null
$d.Content.Find.Execute("32+65=97", $true, $false, $false, $false, $false, $true, 1, $false, "76-1=75", 2) | Out-Null
$d.Content.Find.Execute("11+59=70", $true, $false, $false, $false, $false, $true, 1, $false, "74-19=55", 2) | Out-Null
$d.Content.Find.Execute("18+49=67", $true, $false, $false, $false, $false, $true, 1, $false, "19-5=14", 2) | Out-Null
$d.Content.Find.Execute("63-38=25", $true, $false, $false, $false, $false, $true, 1, $false, "30+66=96", 2) | Out-Null
$d.Content.Find.Execute("38+45=83", $true, $false, $false, $false, $false, $true, 1, $false, "96-5=91", 2) | Out-Null
$d.Content.Find.Execute("29+15=44", $true, $false, $false, $false, $false, $true, 1, $false, "23+22=45", 2) | Out-Null
$d.Content.Find.Execute("6+76=82", $true, $false, $false, $false, $false, $true, 1, $false, "39-0=39", 2) | Out-Null
$d.Content.Find.Execute("58-3=55", $true, $false, $false, $false, $false, $true, 1, $false, "82+10=92", 2) | Out-Null
$d.Content.Find.Execute("79-55=24", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=78", 2) | Out-Null
$d.Content.Find.Execute("43+54=97", $true, $false, $false, $false, $false, $true, 1, $false, "11+11=22", 2) | Out-Null
$d.Content.Find.Execute("40-34=6", $true, $false, $false, $false, $false, $true, 1, $false, "1+45=46", 2) | Out-Null
$d.Content.Find.Execute("99-23=76", $true, $false, $false, $false, $false, $true, 1, $false, "79-60=19", 2) | Out-Null
$d.Content.Find.Execute("68-17=51", $true, $false, $false, $false, $false, $true, 1, $false, "10+21=31", 2) | Out-Null
$d.Content.Find.Execute("72-31=41", $true, $false, $false, $false, $false, $true, 1, $false, "23-9=14", 2) | Out-Null
$d.Content.Find.Execute("21+74=95", $true, $false, $false, $false, $false, $true, 1, $false, "3+50=53", 2) | Out-Null
$d.Content.Find.Execute("74-41=33", $true, $false, $false, $false, $false, $true, 1, $false, "76-39=37", 2) | Out-Null
$d.Content.Find.Execute("4+7=11", $true, $false, $false, $false, $false, $true, 1, $false, "41-40=1", 2) | Out-Null
$d.Content.Find.Execute("35-21=14", $true, $false, $false, $false, $false, $true, 1, $false, "64-33=31", 2) | Out-Null
$d.Content.Find.Execute("7+82=89", $true, $false, $false, $false, $false, $true, 1, $false, "28-0=28", 2) | Out-Null
$d.Content.Find.Execute("82-46=36", $true, $false, $false, $false, $false, $true, 1, $false, "8+31=39", 2) | Out-Null
$d.Content.Find.Execute("73-72=1", $true, $false, $false, $false, $false, $true, 1, $false, "96-40=56", 2) | Out-Null
$d.Content.Find.Execute("66-32=34", $true, $false, $false, $false, $false, $true, 1, $false, "16+53=69", 2) | Out-Null
$d.Content.Find.Execute("22+60=82", $true, $false, $false, $false, $false, $true, 1, $false, "62+6=68", 2) | Out-Null
$d.Content.Find.Execute("31+27=58", $true, $false, $false, $false, $false, $true, 1, $false, "57+39=96", 2) | Out-Null
$d.Content.Find.Execute("45+50=95", $true, $false, $false, $false, $false, $true, 1, $false, "10+1=11", 2) | Out-Null
$d.Content.Find.Execute("20-12=8", $true, $false, $false, $false, $false, $true, 1, $false, "39+31=70", 2) | Out-Null
$d.Content.Find.Execute("85-60=25", $true, $false, $false, $false, $false, $true, 1, $false, "95-92=3", 2) | Out-Null
$d.Content.Find.Execute("44+30=74", $true, $false, $false, $false, $false, $true, 1, $false, "86-41=45", 2) | Out-Null
$d.Content.Find.Execute("27+46=73", $true, $false, $false, $false, $false, $true, 1, $false, "61-57=4", 2) | Out-Null
$d.Content.Find.Execute("48-4=44", $true, $false, $false, $false, $false, $true, 1, $false, "71+28=99", 2) | Out-Null
$d.Content.Find.Execute("4+38=42", $true, $false, $false, $false, $false, $true, 1, $false, "6+69=75", 2) | Out-Null
$d.Content.Find.Execute("32+63=95", $true, $false, $false, $false, $false, $true, 1, $false, "56+38=94", 2) | Out-Null
$d.Content.Find.Execute("19+8=27", $true, $false, $false, $false, $false, $true, 1, $false, "46+24=70", 2) | Out-Null
$d.Content.Find.Execute("39+16=55", $true, $false, $false, $false, $false, $true, 1, $false, "26+36=62", 2) | Out-Null
$d.Content.Find.Execute("61-1=60", $true, $false, $false, $false, $false, $true, 1, $false, "35+41=76", 2) | Out-Null
$d.Content.Find.Execute("42+22=64", $true, $false, $false, $false, $false, $true, 1, $false, "34-12=22", 2) | Out-Null
$d.Content.Find.Execute("57-21=36", $true, $false, $false, $false, $false, $true, 1, $false, "10+20=30", 2) | Out-Null
$d.Content.Find.Execute("30-10=20", $true, $false, $false, $false, $false, $true, 1, $false, "88-2=86", 2) | Out-Null
$d.Content.Find.Execute("67-6=61", $true, $false, $false, $false, $false, $true, 1, $false, "26+69=95", 2) | Out-Null
$d.Content.Find.Execute("23+33=56", $true, $false, $false, $false, $false, $true, 1, $false, "82+0=82", 2) | Out-Null
$d.Content.Find.Execute("48-34=14", $true, $false, $false, $false, $false, $true, 1, $false, "51-19=32", 2) | Out-Null
$d.Content.Find.Execute("72-64=8", $true, $false, $false, $false, $false, $true, 1, $false, "14-10=4", 2) | Out-Null
$d.Content.Find.Execute("6+19=25", $true, $false, $false, $false, $false, $true, 1, $false, "1+7=8", 2) | Out-Null
$d.Content.Find.Execute("45-20=25", $true, $false, $false, $false, $false, $true, 1, $false, "91-42=49", 2) | Out-Null
$d.Content.Find.Execute("7+85=92", $true, $false, $false, $false, $false, $true, 1, $false, "56+19=75", 2) | Out-Null
$d.Content.Find.Execute("76-26=50", $true, $false, $false, $false, $false, $true, 1, $false, "38+44=82", 2) | Out-Null
$d.Content.Find.Execute("85+2=87", $true, $false, $false, $false, $false, $true, 1, $false, "13+66=79", 2) | Out-Null
$d.Content.Find.Execute("79-53=26", $true, $false, $false, $false, $false, $true, 1, $false, "28+40=68", 2) | Out-Null
$d.Content.Find.Execute("74-26=48", $true, $false, $false, $false, $false, $true, 1, $false, "13+25=38", 2) | Out-Null
$d.Content.Find.Execute("81+10=91", $true, $false, $false, $false, $false, $true, 1, $false, "28+4=32", 2) | Out-Null
$d.Content.Find.Execute("92-0=92", $true, $false, $false, $false, $false, $true, 1, $false, "62-9=53", 2) | Out-Null
$d.Content.Find.Execute("35+49=84", $true, $false, $false, $false, $false, $true, 1, $false, "56-51=5", 2) | Out-Null
$d.Content.Find.Execute("92-84=8", $true, $false, $false, $false, $false, $true, 1, $false, "24+52=76", 2) | Out-Null
$d.Content.Find.Execute("45+54=99", $true, $false, $false, $false, $false, $true, 1, $false, "27+33=60", 2) | Out-Null
$d.Content.Find.Execute("3+92=95", $true, $false, $false, $false, $false, $true, 1, $false, "37+8=45", 2) | Out-Null
$d.Content.Find.Execute("92-34=58", $true, $false, $false, $false, $false, $true, 1, $false, "92-45=47", 2) | Out-Null
$d.Content.Find.Execute("20+62=82", $true, $false, $false, $false, $false, $true, 1, $false, "81-18=63", 2) | Out-Null
$d.Content.Find.Execute("90-62=28", $true, $false, $false, $false, $false, $true, 1, $false, "25+10=35", 2) | Out-Null
$d.Content.Find.Execute("45+49=94", $true, $false, $false, $false, $false, $true, 1, $false, "79-35=44", 2) | Out-Null
$d.Content.Find.Execute("42+43=85", $true, $false, $false, $false, $false, $true, 1, $false, "30+5=35", 2) | Out-Null
$d.Content.Find.Execute("32-29=3", $true, $false, $false, $false, $false, $true, 1, $false, "89-4=85", 2) | Out-Null
$d.Content.Find.Execute("33+28=61", $true, $false, $false, $false, $false, $true, 1, $false, "84-52=32", 2) | Out-Null
$d.Content.Find.Execute("21+51=72", $true, $false, $false, $false, $false, $true, 1, $false, "1+78=79", 2) | Out-Null
$d.Content.Find.Execute("9+73=82", $true, $false, $false, $false, $false, $true, 1, $false, "67+9=76", 2) | Out-Null
$d.Content.Find.Execute("50-7=43", $true, $false, $false, $false, $false, $true, 1, $false, "80-38=42", 2) | Out-Null
$d.Content.Find.Execute("10+86=96", $true, $false, $false, $false, $false, $true, 1, $false, "27+48=75", 2) | Out-Null
$d.Content.Find.Execute("95-41=54", $true, $false, $false, $false, $false, $true, 1, $false, "48+17=65", 2) | Out-Null
$d.Content.Find.Execute("8+42=50", $true, $false, $false, $false, $false, $true, 1, $false, "42-41=1", 2) | Out-Null
$d.Content.Find.Execute("24-16=8", $true, $false, $false, $false, $false, $true, 1, $false, "83+11=94", 2) | Out-Null
$d.Content.Find.Execute("77-63=14", $true, $false, $false, $false, $false, $true, 1, $false, "4+95=99", 2) | Out-Null
$d.Content.Find.Execute("2+28=30", $true, $false, $false, $false, $false, $true, 1, $false, "13+24=37", 2) | Out-Null
